## Frassinoro.xlsx report update (commit: "aggiornato a 2/3, aggiornati i report")
##
## The daily case-count table was missing 2021-02-08 (serial 44235); it is inserted here,
## which shifts every later date down one row, and the report window is extended by one
## extra (still-empty) day at the end, 2021-03-02 (serial 44257). Both changes ripple through
## the "somma mobile 7gg." (7-day rolling sum) and "...per 100mila abitanti" columns for the
## dates whose trailing 7-day window now includes the newly inserted day. Column B is the
## per-day new-case count; only 2021-02-08 and the date that used to hold the single new case
## (now one row later) are non-zero. There are no formulas in this sheet -- the report tool
## that produces it bakes every value in -- so we just (re)write the literal numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (dates) uses a dedicated style (border + bold + centered + custom date number
# format) on every existing data row. Clone that formatting onto the two brand-new rows at
# the bottom (114, 115) so their date cells match the rest of the column instead of falling
# back to the default/general style.
$ws.Range("A112").Copy()
$ws.Range("A114:A115").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 90: 44232
$ws.Range("A90").Value = 44232
$ws.Range("B90").Value = 0
$ws.Range("C90").Value = 1
$ws.Range("D90").Value = 53.73455131649651

# Row 91: 44233
$ws.Range("A91").Value = 44233
$ws.Range("B91").Value = 0
$ws.Range("C91").Value = 1
$ws.Range("D91").Value = 53.73455131649651

# Row 92: 44234
$ws.Range("A92").Value = 44234
$ws.Range("B92").Value = 0
$ws.Range("C92").Value = 1
$ws.Range("D92").Value = 53.73455131649651

# Row 93: 44235
$ws.Range("A93").Value = 44235
$ws.Range("B93").Value = 1
$ws.Range("C93").Value = 1
$ws.Range("D93").Value = 53.73455131649651

# Row 94: 44236
$ws.Range("A94").Value = 44236
$ws.Range("B94").Value = 0
$ws.Range("C94").Value = 1
$ws.Range("D94").Value = 53.73455131649651

# Row 95: 44237
$ws.Range("A95").Value = 44237
$ws.Range("B95").Value = 0
$ws.Range("C95").Value = 1
$ws.Range("D95").Value = 53.73455131649651

# Row 96: 44238
$ws.Range("A96").Value = 44238
$ws.Range("B96").Value = 0
$ws.Range("C96").Value = 1
$ws.Range("D96").Value = 53.73455131649651

# Row 97: 44239
$ws.Range("A97").Value = 44239
$ws.Range("B97").Value = 0
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 0

# Row 98: 44240
$ws.Range("A98").Value = 44240
$ws.Range("B98").Value = 0
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 0

# Row 99: 44241
$ws.Range("A99").Value = 44241
$ws.Range("B99").Value = 0
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 0

# Row 100: 44242
$ws.Range("A100").Value = 44242
$ws.Range("B100").Value = 0
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 0

# Row 101: 44243
$ws.Range("A101").Value = 44243
$ws.Range("B101").Value = 0
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 0

# Row 102: 44244
$ws.Range("A102").Value = 44244
$ws.Range("B102").Value = 0
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 0

# Row 103: 44245
$ws.Range("A103").Value = 44245
$ws.Range("B103").Value = 0
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 0

# Row 104: 44246
$ws.Range("A104").Value = 44246
$ws.Range("B104").Value = 0
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 0

# Row 105: 44247
$ws.Range("A105").Value = 44247
$ws.Range("B105").Value = 0
$ws.Range("C105").Value = 1
$ws.Range("D105").Value = 53.73455131649651

# Row 106: 44248
$ws.Range("A106").Value = 44248
$ws.Range("B106").Value = 0
$ws.Range("C106").Value = 1
$ws.Range("D106").Value = 53.73455131649651

# Row 107: 44249
$ws.Range("A107").Value = 44249
$ws.Range("B107").Value = 0
$ws.Range("C107").Value = 1
$ws.Range("D107").Value = 53.73455131649651

# Row 108: 44250
$ws.Range("A108").Value = 44250
$ws.Range("B108").Value = 1
$ws.Range("C108").Value = 1
$ws.Range("D108").Value = 53.73455131649651

# Row 109: 44251
$ws.Range("A109").Value = 44251
$ws.Range("B109").Value = 0
$ws.Range("C109").Value = 1
$ws.Range("D109").Value = 53.73455131649651

# Row 110: 44252
$ws.Range("A110").Value = 44252
$ws.Range("B110").Value = 0
$ws.Range("C110").Value = 1
$ws.Range("D110").Value = 53.73455131649651

# Row 111: 44253
$ws.Range("A111").Value = 44253
$ws.Range("B111").Value = 0
$ws.Range("C111").Value = 1
$ws.Range("D111").Value = 53.73455131649651

# Row 112: 44254
$ws.Range("A112").Value = 44254
$ws.Range("B112").Value = 0
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 0

# Row 113: 44255
$ws.Range("A113").Value = 44255
$ws.Range("B113").Value = 0

# Row 114: 44256
$ws.Range("A114").Value = 44256
$ws.Range("B114").Value = 0

# Row 115: 44257
$ws.Range("A115").Value = 44257
$ws.Range("B115").Value = 0
